$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2279
$ws.Range("I32").Value = 1948.5
$ws.Range("J32").Value = 2361.625
$ws.Range("K32").Value = 1948.5
$ws.Range("L32").Value = 2361.625
$ws.Range("M32").Value = -1622.5
$ws.Range("N32").Value = -3013.625
$ws.Range("H40").Value = 5190.05
$ws.Range("I40").Value = 2650.0833
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 2650.0833
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -2475.0833
$ws.Range("N40").Value = -9350
$ws.Range("H132").Value = 5472.385
$ws.Range("I132").Value = 5470.0835
$ws.Range("K132").Value = 16410.2505
$ws.Range("M132").Value = -13880.2505
$ws.Range("H135").Value = 2077.5386
$ws.Range("I135").Value = 1477.3334
$ws.Range("K135").Value = 13296.0006
$ws.Range("M135").Value = -10761.0006
$ws.Range("H137").Value = 5467.636
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5467.636
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 16402.908
$ws.Range("N137").Value = -21502.908
$ws.Range("H138").Value = 7142.927
$ws.Range("I138").Value = 5179.1665
$ws.Range("J138").Value = 7479.5713
$ws.Range("K138").Value = 15537.4995
$ws.Range("L138").Value = 22438.7139
$ws.Range("M138").Value = -10397.4995
$ws.Range("N138").Value = -32718.7139
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10089.9
$ws.Range("I45").Value = 2966.6667
$ws.Range("J45").Value = 13142.714
$ws.Range("K45").Value = 2966.6667
$ws.Range("L45").Value = 13142.714
$ws.Range("M45").Value = -2589.6667
$ws.Range("N45").Value = -13896.714
$ws.Range("H69").Value = 186000
$ws.Range("J69").Value = 186000
$ws.Range("L69").Value = 186000
$ws.Range("N69").Value = -187498
$ws.Range("H72").Value = 186000
$ws.Range("J72").Value = 186000
$ws.Range("L72").Value = 558000
$ws.Range("N72").Value = -565488
$ws.Range("H122").Value = 3623.6182
$ws.Range("I122").Value = 3612.2693
$ws.Range("K122").Value = 10836.8079
$ws.Range("M122").Value = -8386.8079
$ws.Range("H132").Value = 11358.777
$ws.Range("I132").Value = 4439.8887
$ws.Range("J132").Value = 18277.666
$ws.Range("K132").Value = 13319.6661
$ws.Range("L132").Value = 54832.99800000001
$ws.Range("M132").Value = -10789.6661
$ws.Range("N132").Value = -59892.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1126279.2
$ws.Range("I86").Value = 1636896.8
$ws.Range("K86").Value = 1636896.8
$ws.Range("M86").Value = -1635773.8
$ws.Range("H89").Value = 1126279.2
$ws.Range("I89").Value = 1636896.8
$ws.Range("K89").Value = 8184484
$ws.Range("M89").Value = -8178868
$ws.Range("H134").Value = 40698.793
$ws.Range("I134").Value = 6382.4736
$ws.Range("K134").Value = 19147.4208
$ws.Range("M134").Value = -16612.4208
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3935.5715
$ws.Range("I62").Value = 3931.8
$ws.Range("J62").Value = 3945
$ws.Range("K62").Value = 3931.8
$ws.Range("L62").Value = 3945
$ws.Range("M62").Value = -3307.8
$ws.Range("N62").Value = -5193
$ws.Range("H65").Value = 3935.5715
$ws.Range("I65").Value = 3931.8
$ws.Range("J65").Value = 3945
$ws.Range("K65").Value = 19659
$ws.Range("L65").Value = 19725
$ws.Range("M65").Value = -16539
$ws.Range("N65").Value = -25965
$ws.Range("H86").Value = 5132.1665
$ws.Range("I86").Value = 3948
$ws.Range("J86").Value = 5724.25
$ws.Range("K86").Value = 3948
$ws.Range("L86").Value = 5724.25
$ws.Range("M86").Value = -2825
$ws.Range("N86").Value = -7970.25
$ws.Range("H89").Value = 5132.1665
$ws.Range("I89").Value = 3948
$ws.Range("J89").Value = 5724.25
$ws.Range("K89").Value = 19740
$ws.Range("L89").Value = 28621.25
$ws.Range("M89").Value = -14124
$ws.Range("N89").Value = -39853.25
$ws.Range("H107").Value = 1378.6207
$ws.Range("I107").Value = 692.6667
$ws.Range("J107").Value = 2501.0908
$ws.Range("K107").Value = 692.6667
$ws.Range("L107").Value = 2501.0908
$ws.Range("M107").Value = 1227.3333
$ws.Range("N107").Value = -6341.0908
$ws.Range("H132").Value = 3574.9546
$ws.Range("I132").Value = 3332.5
$ws.Range("J132").Value = 5999.5
$ws.Range("K132").Value = 9997.5
$ws.Range("L132").Value = 17998.5
$ws.Range("M132").Value = -7467.5
$ws.Range("N132").Value = -23058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 93.7
$ws.Range("I10").Value = 61.625
$ws.Range("J10").Value = 222
$ws.Range("K10").Value = 184.875
$ws.Range("L10").Value = 666
$ws.Range("M10").Value = -45.875
$ws.Range("N10").Value = -944
$ws.Range("H12").Value = 232.5
$ws.Range("I12").Value = 12.5
$ws.Range("K12").Value = 37.5
$ws.Range("M12").Value = 135.5
$ws.Range("H107").Value = 5866.3
$ws.Range("I107").Value = 805
$ws.Range("K107").Value = 2415
$ws.Range("M107").Value = -495
$ws.Range("H122").Value = 43614.543
$ws.Range("J122").Value = 54889.844
$ws.Range("L122").Value = 494008.596
$ws.Range("N122").Value = -498908.596

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 921983.6
$ws.Range("J113").Value = 4973.2856
$ws.Range("L113").Value = 4973.2856
$ws.Range("N113").Value = -9313.285599999999
$ws.Range("H126").Value = 5651.1
$ws.Range("J126").Value = 5833.222
$ws.Range("L126").Value = 17499.666
$ws.Range("N126").Value = -22439.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 672890.5600000001
$ws.Range("I7").Value = 1116262
$ws.Range("K7").Value = 1116262
$ws.Range("M7").Value = -1116150
$ws.Range("H20").Value = 47142.285
$ws.Range("J20").Value = 47142.285
$ws.Range("L20").Value = 47142.285
$ws.Range("N20").Value = -47594.285
$ws.Range("H22").Value = 911.9375
$ws.Range("J22").Value = 999.4
$ws.Range("L22").Value = 999.4
$ws.Range("N22").Value = -1589.4
$ws.Range("H27").Value = 911.9375
$ws.Range("J27").Value = 999.4
$ws.Range("L27").Value = 999.4
$ws.Range("N27").Value = -1213.4
$ws.Range("H40").Value = 775060
$ws.Range("I40").Value = 1115531.2
$ws.Range("K40").Value = 1115531.2
$ws.Range("M40").Value = -1115395.2
$ws.Range("H61").Value = 5793.2856
$ws.Range("I61").Value = 4505.9165
$ws.Range("K61").Value = 4505.9165
$ws.Range("M61").Value = -4303.9165
$ws.Range("H113").Value = 5793.2856
$ws.Range("I113").Value = 4505.9165
$ws.Range("K113").Value = 4505.9165
$ws.Range("M113").Value = -2335.9165
$ws.Range("H126").Value = 672890.5600000001
$ws.Range("I126").Value = 1116262
$ws.Range("K126").Value = 3348786
$ws.Range("M126").Value = -3346316
$ws.Range("H132").Value = 6584.25
$ws.Range("I132").Value = 4732.3335
$ws.Range("J132").Value = 7695.4
$ws.Range("K132").Value = 14197.0005
$ws.Range("L132").Value = 23086.2
$ws.Range("M132").Value = -11667.0005
$ws.Range("N132").Value = -28146.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 815.43243
$ws.Range("I113").Value = 714
$ws.Range("K113").Value = 2142
$ws.Range("M113").Value = 28
$ws.Range("H122").Value = 33338064
$ws.Range("I122").Value = 55558940
$ws.Range("K122").Value = 166676820
$ws.Range("M122").Value = -166674370
$ws.Range("H126").Value = 3014.25
$ws.Range("I126").Value = 1408.2307
$ws.Range("J126").Value = 5996.857
$ws.Range("K126").Value = 4224.6921
$ws.Range("L126").Value = 17990.571
$ws.Range("M126").Value = -1754.6921
$ws.Range("N126").Value = -22930.571
$ws.Range("H132").Value = 52452.43
$ws.Range("I132").Value = 3408.4
$ws.Range("J132").Value = 97037.91
$ws.Range("K132").Value = 10225.2
$ws.Range("L132").Value = 291113.73
$ws.Range("M132").Value = -7695.200000000001
$ws.Range("N132").Value = -296173.73
$ws.Range("H136").Value = 484450.56
$ws.Range("I136").Value = 669764.3
$ws.Range("J136").Value = 206479.9
$ws.Range("K136").Value = 2009292.9
$ws.Range("L136").Value = 619439.7
$ws.Range("M136").Value = -2006742.9
$ws.Range("N136").Value = -624539.7
